$d = $word.ActiveDocument

# 1) Title paragraph: remove indentation (left=0, firstLine=0)
$p1 = $d.Paragraphs.Item(1)
$p1.Format.LeftIndent = 0
$p1.Format.FirstLineIndent = 0

# 2) "O problema ..." paragraph: justify + replace two text runs
$p3 = $d.Paragraphs.Item(3)
$p3.Format.Alignment = 3

$d.Content.Find.Execute(
    "da dificuldade de fazer a gestão dos produtos e outras áreas da empresa ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "da dificuldade de fazer a gestão de produtos e compras da empresa ", 2) | Out-Null

$d.Content.Find.Execute(
    "não existir um site para melhor comunicação e exibição dos produtos da loja, tendo algumas reclamações por parte dos clientes por não ter um sistema próprio de rastreamento dos pedidos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a perda de vendas e poucas avaliações de seus produtos da loja e o cliente não consegue ser atendido imediatamente, tendo que esperar mais tempo para conseguir resposta sobre os produtos e assim, perdendo alguns clientes.", 2) | Out-Null

# 3) empty spacer paragraph after -> justify
$p4 = $d.Paragraphs.Item(4)
$p4.Format.Alignment = 3

# 4) "Os benefícios..." paragraph -> justify
$p5 = $d.Paragraphs.Item(5)
$p5.Format.Alignment = 3

# 5) list items 6,7,8 -> justify only
$p6 = $d.Paragraphs.Item(6)
$p6.Format.Alignment = 3

$p7 = $d.Paragraphs.Item(7)
$p7.Format.Alignment = 3

$p8 = $d.Paragraphs.Item(8)
$p8.Format.Alignment = 3

# 6) list item 9 "Melhor visibilidade da loja;" -> justify + replace text
$p9 = $d.Paragraphs.Item(9)
$p9.Format.Alignment = 3
$d.Content.Find.Execute(
    "Melhor visibilidade da loja;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Facilidade para que o cliente possa cadastrar seus dados no site para futuros pedidos;", 2) | Out-Null

# 7) list item 10 (old "Facilidade para que o cliente possa cadastrar...") -> delete whole paragraph
$p10 = $d.Paragraphs.Item(10)
$p10.Range.Delete()

# 8) list item (now shifted to 10, was 11) "Melhorar a logística..." -> justify only
$p10b = $d.Paragraphs.Item(10)
$p10b.Format.Alignment = 3

# 9) final paragraph (now 11, was 12) -> justify
$pLast = $d.Paragraphs.Item(11)
$pLast.Format.Alignment = 3

Write-Output "Done"
